$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4) to the sheet, matching the existing data pattern.
$ws.Cells.Item(4, 1).Value = 42633.679085648146
$ws.Cells.Item(4, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(4, 2).Value = $false

$ws.Cells.Item(4, 3).Value = 9974
$ws.Cells.Item(4, 4).Value = 10000
$ws.Cells.Item(4, 5).Value = 19.32
$ws.Cells.Item(4, 6).Value = 19.22

$ws.Cells.Item(4, 7).Value = $false
$ws.Cells.Item(4, 7).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(4, 8).Value = -0.52

$ws.Cells.Item(4, 9).Value = $false
